$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Update the existing "PLS_predict" row (row 5) coefficients.
#    These cells hold their numeric-looking values as TEXT (apostrophe
#    prefix). Cells C5/K5 have no special number format (style index 0),
#    while D5/F5/G5/H5/I5/J5 use the workbook's "quote prefix" text style.
#    Writing with a leading apostrophe preserves/produces that quote-prefix
#    style; C5/K5 are then reset back to the "Normal" style so they stay
#    styleless, matching their original formatting.
# ---------------------------------------------------------------------------
$ws.Range("C5").Value = "'0.428"
$ws.Range("C5").Style = "Normal"

$ws.Range("D5").Value = "'-0.433"
$ws.Range("F5").Value = "'0.367"
$ws.Range("G5").Value = "'-0.143"
$ws.Range("H5").Value = "'0.27"
$ws.Range("I5").Value = "'-0.128"
$ws.Range("J5").Value = "'-0.181"

$ws.Range("K5").Value = "'0.593"
$ws.Range("K5").Style = "Normal"

# ---------------------------------------------------------------------------
# 2) Add the new "PLS LOG" model row (row 6).
# ---------------------------------------------------------------------------
$ws.Range("A6").Value = "PLS LOG"
$ws.Range("B6").Value = "PLS with log of turb"

$ws.Range("C6").Value = "'0.389"
$ws.Range("C6").Style = "Normal"

$ws.Range("D6").Value = "'-0.446"
$ws.Range("E6").Value = "'-0.259"
$ws.Range("F6").Value = "'0.300"
$ws.Range("G6").Value = "'-0.110"
$ws.Range("H6").Value = "'0.256"
$ws.Range("I6").Value = "'-0.176"
$ws.Range("J6").Value = "'-0.243"

$ws.Range("K6").Value = "'0.570"
$ws.Range("K6").Style = "Normal"

# ---------------------------------------------------------------------------
# 3) Move the active selection, mirroring the saved view state in the diff.
# ---------------------------------------------------------------------------
$ws.Range("D9").Select() | Out-Null
